{"js": "// Replace each three-digit-times-one-digit multiplication answer\n// with its new value. Every source string in this table is unique,\n// so a case-sensitive exact search safely targets exactly one cell.\nconst replacements = [\n  [\"832\u00d75=4160\", \"160\u00d78=1280\"],\n  [\"556\u00d74=2224\", \"794\u00d77=5558\"],\n  [\"231\u00d79=2079\", \"733\u00d78=5864\"],\n  [\"462\u00d73=1386\", \"858\u00d72=1716\"],\n  [\"552\u00d74=2208\", \"324\u00d77=2268\"],\n  [\"548\u00d76=3288\", \"950\u00d77=6650\"],\n  [\"943\u00d77=6601\", \"325\u00d76=1950\"],\n  [\"817\u00d76=4902\", \"405\u00d78=3240\"],\n  [\"577\u00d78=4616\", \"241\u00d79=2169\"],\n  [\"182\u00d73=546\", \"287\u00d77=2009\"],\n  [\"107\u00d79=963\", \"986\u00d76=5916\"],\n  [\"899\u00d72=1798\", \"195\u00d78=1560\"],\n  [\"301\u00d73=903\", \"393\u00d75=1965\"],\n  [\"453\u00d73=1359\", \"478\u00d72=956\"],\n  [\"834\u00d74=3336\", \"808\u00d74=3232\"],\n  [\"534\u00d77=3738\", \"937\u00d75=4685\"],\n  [\"249\u00d73=747\", \"293\u00d73=879\"],\n  [\"158\u00d75=790\", \"630\u00d74=2520\"],\n  [\"765\u00d79=6885\", \"998\u00d78=7984\"],\n  [\"870\u00d76=5220\", \"388\u00d73=1164\"],\n  [\"505\u00d79=4545\", \"519\u00d79=4671\"],\n  [\"993\u00d75=4965\", \"783\u00d78=6264\"],\n  [\"927\u00d79=8343\", \"404\u00d73=1212\"],\n  [\"368\u00d78=2944\", \"199\u00d78=1592\"],\n  [\"389\u00d78=3112\", \"124\u00d73=372\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update the three-digit-times-one-digit multiplication answers in the\n# practice table. Every source string is unique in the document, so a\n# simple Find/Replace per pair safely targets exactly one cell each.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"832\u00d75=4160\", \"160\u00d78=1280\"),\n    @(\"556\u00d74=2224\", \"794\u00d77=5558\"),\n    @(\"231\u00d79=2079\", \"733\u00d78=5864\"),\n    @(\"462\u00d73=1386\", \"858\u00d72=1716\"),\n    @(\"552\u00d74=2208\", \"324\u00d77=2268\"),\n    @(\"548\u00d76=3288\", \"950\u00d77=6650\"),\n    @(\"943\u00d77=6601\", \"325\u00d76=1950\"),\n    @(\"817\u00d76=4902\", \"405\u00d78=3240\"),\n    @(\"577\u00d78=4616\", \"241\u00d79=2169\"),\n    @(\"182\u00d73=546\", \"287\u00d77=2009\"),\n    @(\"107\u00d79=963\", \"986\u00d76=5916\"),\n    @(\"899\u00d72=1798\", \"195\u00d78=1560\"),\n    @(\"301\u00d73=903\", \"393\u00d75=1965\"),\n    @(\"453\u00d73=1359\", \"478\u00d72=956\"),\n    @(\"834\u00d74=3336\", \"808\u00d74=3232\"),\n    @(\"534\u00d77=3738\", \"937\u00d75=4685\"),\n    @(\"249\u00d73=747\", \"293\u00d73=879\"),\n    @(\"158\u00d75=790\", \"630\u00d74=2520\"),\n    @(\"765\u00d79=6885\", \"998\u00d78=7984\"),\n    @(\"870\u00d76=5220\", \"388\u00d73=1164\"),\n    @(\"505\u00d79=4545\", \"519\u00d79=4671\"),\n    @(\"993\u00d75=4965\", \"783\u00d78=6264\"),\n    @(\"927\u00d79=8343\", \"404\u00d73=1212\"),\n    @(\"368\u00d78=2944\", \"199\u00d78=1592\"),\n    @(\"389\u00d78=3112\", \"124\u00d73=372\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n    if (-not $found) {\n        throw \"Search text not found: $old\"\n    }\n}\n\n"}
